# Update dashboards - 2025-10-18
# Applies the daily data refresh to the "Aguilar Prototype" sheet:
#  - Row 29 (T5YIFR, monthly): latest date revised, Present value revised
#  - Row 30 (T10YIE, daily): date advances one day, 5-day window shifts
#  - Row 48 (DGS2, daily): date advances one day, 5-day window shifts
#  - Row 49 (DGS5, daily): date advances one day, 5-day window shifts
#  - Row 50 (DGS10, daily): date advances one day, 5-day window shifts
#  - Row 52 (DBAA, daily): date advances one day, 5-day window shifts

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DateText {
    param($addr, $text)
    # Prefix with an apostrophe so Excel stores the value as text rather than
    # auto-converting the yyyy-mm-dd looking string into a date serial number.
    $ws.Range($addr).Value = "'" + $text
}

# --- Row 29 : 5yr,5yr Forward (T5YIFR) - monthly series, value revision only ---
Set-DateText "N29" "2025-10-17"
$ws.Range("Q29").Value = 2.25

# --- Row 30 : 10yr TIPS (T10YIE) - daily series, window shift ---
Set-DateText "N30" "2025-10-17"
$ws.Range("Q30").Value = 2.27
$ws.Range("R30").Value = 2.28
$ws.Range("S30").Value = 2.29
$ws.Range("T30").Value = 2.3
$ws.Range("U30").ClearContents()

# --- Row 48 : 2y UST (DGS2) - daily series, window shift ---
Set-DateText "N48" "2025-10-16"
$ws.Range("Q48").Value = 3.41
$ws.Range("R48").Value = 3.5
$ws.Range("S48").Value = 3.48
$ws.Range("T48").ClearContents()
$ws.Range("U48").Value = 3.52

# --- Row 49 : 5y UST (DGS5) - daily series, window shift ---
Set-DateText "N49" "2025-10-16"
$ws.Range("Q49").Value = 3.55
$ws.Range("R49").Value = 3.63
$ws.Range("S49").Value = 3.6
$ws.Range("T49").ClearContents()
$ws.Range("U49").Value = 3.65

# --- Row 50 : 10y UST (DGS10) - daily series, window shift ---
Set-DateText "N50" "2025-10-16"
$ws.Range("Q50").Value = 3.99
$ws.Range("R50").Value = 4.05
$ws.Range("S50").Value = 4.03
$ws.Range("T50").ClearContents()
$ws.Range("U50").Value = 4.05

# --- Row 52 : BAA (DBAA) - daily series, window shift ---
Set-DateText "N52" "2025-10-16"
$ws.Range("Q52").Value = 5.7
$ws.Range("R52").Value = 5.73
$ws.Range("S52").Value = 5.74
$ws.Range("T52").ClearContents()
$ws.Range("U52").Value = 5.77
